# Apply the Saldo.xlsx update:
#  - Row 2 (account 005547703 / SILVIA): balance 130000 -> 78861.27
#  - Row 7 (account 004866753 / GENESI / 8213.88) replaced with
#    account 005040864 / ANDRE / 12233.2
#  - Remove the duplicate/old row further down that held
#    005040864 / ANDRE / 100 (row 109)
#  - Remove row for account 001761119 / BLUEMETRIX / 126.98 (row 99)
#  - Remove row for account 004384167 / DOUGLAS / 19294.8 (row 5)
#  - Remove row for account 002694089 / VITOR / 25473.87 (row 3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update SILVIA's balance in-place.
$ws.Range("C2").Value = 78861.27

# 2) Turn the GENESI row into the ANDRE row (new balance) in-place.
#    Force column A to text format first so the leading zeros in the
#    account number are preserved instead of being parsed as a number.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "005040864"
$ws.Range("B7").Value = "ANDRE"
$ws.Range("C7").Value = 12233.2

# 3) Delete rows bottom-up so earlier row numbers stay valid.
$ws.Rows(109).Delete()   # old 005040864 / ANDRE / 100 row
$ws.Rows(99).Delete()    # 001761119 / BLUEMETRIX / 126.98
$ws.Rows(5).Delete()     # 004384167 / DOUGLAS / 19294.8
$ws.Rows(3).Delete()     # 002694089 / VITOR / 25473.87
